# "added range argument to read_excel"
#
# Adds a new worksheet "position" at the end of the workbook containing the
# same data as the "3d" sheet, but written starting at D3 instead of A1
# (used to exercise a `range=` argument on read_excel). Also tweaks a few
# pre-existing sheet selections to match the saved workbook state.

$wb = $excel.ActiveWorkbook

# --- "3d" sheet: selection becomes the whole used range (A1:E7) ----------
$ws3d = $wb.Worksheets.Item("3d")
$ws3d.Activate()
[void]$ws3d.UsedRange.Select()

# --- "int_labels" sheet: new selection at H15 -----------------------------
$wsIntLabels = $wb.Worksheets.Item("int_labels")
$wsIntLabels.Activate()
[void]$wsIntLabels.Range("H15").Select()

# --- new "position" sheet, appended after "int_labels" --------------------
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "position"

$data = @(
    @("a", "b\c", "c0", "c1", "c2"),
    @(1, "b0", 0, 1, 2),
    @(1, "b1", 3, 4, 5),
    @(2, "b0", 6, 7, 8),
    @(2, "b1", 9, 10, 11),
    @(3, "b0", 12, 13, 14),
    @(3, "b1", 15, 16, 17)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $newSheet.Cells.Item(3 + $i, 4 + $j).Value = $row[$j]
    }
}

# Leave the new sheet active, with its selection on the bottom-right cell of
# the data (mirrors how Excel leaves the cursor after typing in the range).
[void]$newSheet.Range("H9").Select()
